$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row labels (A1:D1)
$ws.Range('A1').Value = 'mx_state'
$ws.Range('B1').Value = 'mx_municipality'
$ws.Range('C1').Value = 'n_matriculas'
$ws.Range('D1').Value = 'pct_matriculas'

# Title-case the Spanish connector words ("de", "del", "la", "los", "las", "el", "y")
# in state/municipality names, and fix the final TOTAL label casing.
$ws.Cells.Item(6, 2).Value = 'Pabellón De Arteaga'
$ws.Cells.Item(7, 2).Value = 'Rincón De Romos'
$ws.Cells.Item(8, 2).Value = 'San Francisco De Los Romo'
$ws.Cells.Item(32, 2).Value = 'Amatenango De La Frontera'
$ws.Cells.Item(33, 2).Value = 'Amatenango Del Valle'
$ws.Cells.Item(37, 2).Value = 'Bejucal De Ocampo'
$ws.Cells.Item(39, 2).Value = 'Benemérito De Las Américas'
$ws.Cells.Item(49, 2).Value = 'Chiapa De Corzo'
$ws.Cells.Item(54, 2).Value = 'Comitán De Domínguez'
$ws.Cells.Item(83, 2).Value = 'Marqués De Comillas'
$ws.Cells.Item(84, 2).Value = 'Mazapa De Madero'
$ws.Cells.Item(88, 2).Value = 'Montecristo De Guerrero'
$ws.Cells.Item(92, 2).Value = 'Ocozocoautla De Espinosa'
$ws.Cells.Item(103, 2).Value = 'Salto De Agua'
$ws.Cells.Item(104, 2).Value = 'San Cristóbal De Las Casas'
$ws.Cells.Item(106, 2).Value = 'Santiago El Pinar'
$ws.Cells.Item(146, 2).Value = 'Guadalupe Y Calvo'
$ws.Cells.Item(148, 2).Value = 'Hidalgo Del Parral'
$ws.Cells.Item(156, 1).Value = 'Ciudad De México'
$ws.Cells.Item(159, 2).Value = 'Cuajimalpa De Morelos'
$ws.Cells.Item(174, 1).Value = 'Coahuila De Zaragoza'
$ws.Cells.Item(191, 2).Value = 'San Juan De Sabinas'
$ws.Cells.Item(202, 2).Value = 'Villa De Álvarez'
$ws.Cells.Item(215, 2).Value = 'Nombre De Dios'
$ws.Cells.Item(219, 2).Value = 'Pánuco De Coronado'
$ws.Cells.Item(221, 2).Value = 'San Juan De Guadalupe'
$ws.Cells.Item(222, 2).Value = 'San Juan Del Río'
$ws.Cells.Item(223, 2).Value = 'San Luis Del Cordero'
$ws.Cells.Item(232, 1).Value = 'Estado De México'
$ws.Cells.Item(232, 2).Value = 'Acambay De Ruíz Castañeda'
$ws.Cells.Item(235, 2).Value = 'Almoloya De Alquisiras'
$ws.Cells.Item(236, 2).Value = 'Almoloya De Juárez'
$ws.Cells.Item(242, 2).Value = 'Atizapán De Zaragoza'
$ws.Cells.Item(248, 2).Value = 'Chapa De Mota'
$ws.Cells.Item(251, 2).Value = 'Coacalco De Berriozábal'
$ws.Cells.Item(258, 2).Value = 'Ecatepec De Morelos'
$ws.Cells.Item(265, 2).Value = 'Ixtapan De La Sal'
$ws.Cells.Item(266, 2).Value = 'Ixtapan Del Oro'
$ws.Cells.Item(280, 2).Value = 'Naucalpan De Juárez'
$ws.Cells.Item(290, 2).Value = 'San Felipe Del Progreso'
$ws.Cells.Item(291, 2).Value = 'San José Del Rincón'
$ws.Cells.Item(292, 2).Value = 'San Martín De Las Pirámides'
$ws.Cells.Item(294, 2).Value = 'San Simón De Guerrero'
$ws.Cells.Item(296, 2).Value = 'Soyaniquilpan De Juárez'
$ws.Cells.Item(305, 2).Value = 'Tenango Del Valle'
$ws.Cells.Item(315, 2).Value = 'Tlalnepantla De Baz'
$ws.Cells.Item(320, 2).Value = 'Valle De Bravo'
$ws.Cells.Item(321, 2).Value = 'Valle De Chalco Solidaridad'
$ws.Cells.Item(324, 2).Value = 'Villa De Allende'
$ws.Cells.Item(325, 2).Value = 'Villa Del Carbón'
$ws.Cells.Item(335, 2).Value = 'Apaseo El Alto'
$ws.Cells.Item(336, 2).Value = 'Apaseo El Grande'
$ws.Cells.Item(345, 2).Value = 'Dolores Hidalgo Cuna De La Independencia Nacional'
$ws.Cells.Item(349, 2).Value = 'Jaral Del Progreso'
$ws.Cells.Item(356, 2).Value = 'Purísima Del Rincón'
$ws.Cells.Item(361, 2).Value = 'San Diego De La Unión'
$ws.Cells.Item(363, 2).Value = 'San Francisco Del Rincón'
$ws.Cells.Item(365, 2).Value = 'San Luis De La Paz'
$ws.Cells.Item(366, 2).Value = 'San Miguel De Allende'
$ws.Cells.Item(367, 2).Value = 'Santa Cruz De Juventino Rosas'
$ws.Cells.Item(369, 2).Value = 'Silao De La Victoria'
$ws.Cells.Item(374, 2).Value = 'Valle De Santiago'
$ws.Cells.Item(380, 2).Value = 'Acapulco De Juárez'
$ws.Cells.Item(383, 2).Value = 'Ajuchitlán Del Progreso'
$ws.Cells.Item(384, 2).Value = 'Alcozauca De Guerrero'
$ws.Cells.Item(388, 2).Value = 'Atenango Del Río'
$ws.Cells.Item(389, 2).Value = 'Atlamajalcingo Del Monte'
$ws.Cells.Item(391, 2).Value = 'Atoyac De Álvarez'
$ws.Cells.Item(392, 2).Value = 'Ayutla De Los Libres'
$ws.Cells.Item(395, 2).Value = 'Buenavista De Cuéllar'
$ws.Cells.Item(396, 2).Value = 'Chilapa De Álvarez'
$ws.Cells.Item(397, 2).Value = 'Chilpancingo De Los Bravo'
$ws.Cells.Item(398, 2).Value = 'Coahuayutla De José María Izazaga'
$ws.Cells.Item(403, 2).Value = 'Coyuca De Benítez'
$ws.Cells.Item(404, 2).Value = 'Coyuca De Catalán'
$ws.Cells.Item(408, 2).Value = 'Cuetzala Del Progreso'
$ws.Cells.Item(409, 2).Value = 'Cutzamala De Pinzón'
$ws.Cells.Item(416, 2).Value = 'Huitzuco De Los Figueroa'
$ws.Cells.Item(417, 2).Value = 'Iguala De La Independencia'
$ws.Cells.Item(419, 2).Value = 'Ixcateopan De Cuauhtémoc'
$ws.Cells.Item(422, 2).Value = 'La Unión De Isidoro Montes De Oca'
$ws.Cells.Item(427, 2).Value = 'Mártir De Cuilapan'
$ws.Cells.Item(438, 2).Value = 'Taxco De Alarcón'
$ws.Cells.Item(441, 2).Value = 'Tepecoacuilco De Trujano'
$ws.Cells.Item(443, 2).Value = 'Tixtla De Guerrero'
$ws.Cells.Item(447, 2).Value = 'Tlalixtaquilla De Maldonado'
$ws.Cells.Item(448, 2).Value = 'Tlapa De Comonfort'
$ws.Cells.Item(450, 2).Value = 'Técpan De Galeana'
$ws.Cells.Item(455, 2).Value = 'Zihuatanejo De Azueta'
$ws.Cells.Item(462, 2).Value = 'Agua Blanca De Iturbide'
$ws.Cells.Item(468, 2).Value = 'Atotonilco El Grande'
$ws.Cells.Item(474, 2).Value = 'Cuautepec De Hinojosa'
$ws.Cells.Item(480, 2).Value = 'Huasca De Ocampo'
$ws.Cells.Item(485, 2).Value = 'Huejutla De Reyes'
$ws.Cells.Item(489, 2).Value = 'Jacala De Ledezma'
$ws.Cells.Item(495, 2).Value = 'Mineral De La Reforma'
$ws.Cells.Item(496, 2).Value = 'Mineral Del Chico'
$ws.Cells.Item(497, 2).Value = 'Mineral Del Monte'
$ws.Cells.Item(498, 2).Value = 'Mixquiahuala De Juárez'
$ws.Cells.Item(499, 2).Value = 'Molango De Escamilla'
$ws.Cells.Item(501, 2).Value = 'Nopala De Villagrán'
$ws.Cells.Item(502, 2).Value = 'Omitlán De Juárez'
$ws.Cells.Item(503, 2).Value = 'Pachuca De Soto'
$ws.Cells.Item(506, 2).Value = 'Progreso De Obregón'
$ws.Cells.Item(511, 2).Value = 'Santiago Tulantepec De Lugo Guerrero'
$ws.Cells.Item(512, 2).Value = 'Santiago De Anaya'
$ws.Cells.Item(516, 2).Value = 'Tenango De Doria'
$ws.Cells.Item(518, 2).Value = 'Tepehuacán De Guerrero'
$ws.Cells.Item(519, 2).Value = 'Tepeji Del Río De Ocampo'
$ws.Cells.Item(521, 2).Value = 'Tezontepec De Aldama'
$ws.Cells.Item(529, 2).Value = 'Tula De Allende'
$ws.Cells.Item(530, 2).Value = 'Tulancingo De Bravo'
$ws.Cells.Item(531, 2).Value = 'Villa De Tezontepec'
$ws.Cells.Item(535, 2).Value = 'Zacualtipán De Ángeles'
$ws.Cells.Item(540, 2).Value = 'Acatlán De Juárez'
$ws.Cells.Item(541, 2).Value = 'Ahualulco De Mercado'
$ws.Cells.Item(545, 2).Value = 'Atotonilco El Alto'
$ws.Cells.Item(547, 2).Value = 'Autlán De Navarro'
$ws.Cells.Item(555, 2).Value = 'Cuautitlán De García Barragán'
$ws.Cells.Item(563, 2).Value = 'Encarnación De Díaz'
$ws.Cells.Item(567, 2).Value = 'Huejuquilla El Alto'
$ws.Cells.Item(569, 2).Value = 'Ixtlahuacán De Los Membrillos'
$ws.Cells.Item(570, 2).Value = 'Ixtlahuacán Del Río'
$ws.Cells.Item(573, 2).Value = 'Jilotlán De Los Dolores'
$ws.Cells.Item(577, 2).Value = 'Lagos De Moreno'
$ws.Cells.Item(582, 2).Value = 'Ojuelos De Jalisco'
$ws.Cells.Item(587, 2).Value = 'San Cristóbal De La Barranca'
$ws.Cells.Item(588, 2).Value = 'San Diego De Alejandría'
$ws.Cells.Item(591, 2).Value = 'San Martín De Bolaños'
$ws.Cells.Item(592, 2).Value = 'San Miguel El Alto'
$ws.Cells.Item(594, 2).Value = 'Santa María De Los Ángeles'
$ws.Cells.Item(595, 2).Value = 'Santa María Del Oro'
$ws.Cells.Item(597, 2).Value = 'Tamazula De Gordiano'
$ws.Cells.Item(601, 2).Value = 'Tepatitlán De Morelos'
$ws.Cells.Item(604, 2).Value = 'Tizapán El Alto'
$ws.Cells.Item(605, 2).Value = 'Tlajomulco De Zúñiga'
$ws.Cells.Item(611, 2).Value = 'Unión De San Antonio'
$ws.Cells.Item(612, 2).Value = 'Unión De Tula'
$ws.Cells.Item(613, 2).Value = 'Valle De Juárez'
$ws.Cells.Item(616, 2).Value = 'Yahualica De González Gallo'
$ws.Cells.Item(619, 2).Value = 'Zapotitlán De Vadillo'
$ws.Cells.Item(620, 2).Value = 'Zapotlán Del Rey'
$ws.Cells.Item(621, 2).Value = 'Zapotlán El Grande'
$ws.Cells.Item(623, 1).Value = 'Michoacán De Ocampo'
$ws.Cells.Item(642, 2).Value = 'Coalcomán De Vázquez Pallares'
$ws.Cells.Item(644, 2).Value = 'Cojumatlán De Régules'
$ws.Cells.Item(707, 2).Value = 'Tiquicheo De Nicolás Romero'
$ws.Cells.Item(731, 2).Value = 'Coatlán Del Río'
$ws.Cells.Item(744, 2).Value = 'Puente De Ixtla'
$ws.Cells.Item(749, 2).Value = 'Tetela Del Volcán'
$ws.Cells.Item(751, 2).Value = 'Tlaltizapán De Zapata'
$ws.Cells.Item(762, 2).Value = 'Bahía De Banderas'
$ws.Cells.Item(766, 2).Value = 'Ixtlán Del Río'
$ws.Cells.Item(773, 2).Value = 'Santa María Del Oro'
$ws.Cells.Item(791, 2).Value = 'Lampazos De Naranjo'
$ws.Cells.Item(794, 2).Value = 'Mier Y Noriega'
$ws.Cells.Item(798, 2).Value = 'San Nicolás De Los Garza'
$ws.Cells.Item(802, 2).Value = 'Acatlán De Pérez Figueroa'
$ws.Cells.Item(813, 2).Value = 'El Barrio De La Soledad'
$ws.Cells.Item(816, 2).Value = 'Heroica Ciudad De Ejutla De Crespo'
$ws.Cells.Item(817, 2).Value = 'Heroica Ciudad De Huajuapan De León'
$ws.Cells.Item(818, 2).Value = 'Heroica Ciudad De Juchitán De Zaragoza'
$ws.Cells.Item(819, 2).Value = 'Heroica Ciudad De Tlaxiaco'
$ws.Cells.Item(820, 2).Value = 'Huajuapan De León'
$ws.Cells.Item(822, 2).Value = 'Huautla De Jiménez'
$ws.Cells.Item(823, 2).Value = 'Ixtlán De Juárez'
$ws.Cells.Item(828, 2).Value = 'Mariscala De Juárez'
$ws.Cells.Item(830, 2).Value = 'Mazatlán Villa De Flores'
$ws.Cells.Item(832, 2).Value = 'Miahuatlán De Porfirio Díaz'
$ws.Cells.Item(834, 2).Value = 'Mártires De Tacubaya'
$ws.Cells.Item(835, 2).Value = 'Nejapa De Madero'
$ws.Cells.Item(836, 2).Value = 'Oaxaca De Juárez'
$ws.Cells.Item(837, 2).Value = 'Ocotlán De Morelos'
$ws.Cells.Item(838, 2).Value = 'Pinotepa De Don Luis'
$ws.Cells.Item(840, 2).Value = 'Putla Villa De Guerrero'
$ws.Cells.Item(846, 2).Value = 'San Agustín De Las Juntas'
$ws.Cells.Item(853, 2).Value = 'San Antonino El Alto'
$ws.Cells.Item(856, 2).Value = 'San Antonio De La Cal'
$ws.Cells.Item(861, 2).Value = 'San Dionisio Del Mar'
$ws.Cells.Item(863, 2).Value = 'San Felipe Jalapa De Díaz'
$ws.Cells.Item(869, 2).Value = 'San Francisco Del Mar'
$ws.Cells.Item(879, 2).Value = 'San José Del Progreso'
$ws.Cells.Item(881, 2).Value = 'San Juan Bautista Lo De Soto'
$ws.Cells.Item(913, 2).Value = 'San Mateo Del Mar'
$ws.Cells.Item(924, 2).Value = 'San Miguel Del Puerto'
$ws.Cells.Item(929, 2).Value = 'San Pablo Villa De Mitla'
$ws.Cells.Item(946, 2).Value = 'San Pedro El Alto'
$ws.Cells.Item(947, 2).Value = 'San Pedro Y San Pablo Ayutla'
$ws.Cells.Item(948, 2).Value = 'San Pedro Y San Pablo Teposcolula'
$ws.Cells.Item(973, 2).Value = 'Santa Inés De Zaragoza'
$ws.Cells.Item(974, 2).Value = 'Santa Inés Del Monte'
$ws.Cells.Item(986, 2).Value = 'Santa María Jalapa Del Marqués'
$ws.Cells.Item(1035, 2).Value = 'Santo Domingo De Morelos'
$ws.Cells.Item(1040, 2).Value = 'Sitio De Xitlapehua'
$ws.Cells.Item(1041, 2).Value = 'Tamazulápam Del Espíritu Santo'
$ws.Cells.Item(1043, 2).Value = 'Tataltepec De Valdés'
$ws.Cells.Item(1044, 2).Value = 'Teotitlán De Flores Magón'
$ws.Cells.Item(1046, 2).Value = 'Tezoatlán De Segura Y Luna'
$ws.Cells.Item(1047, 2).Value = 'Tlacolula De Matamoros'
$ws.Cells.Item(1048, 2).Value = 'Totontepec Villa De Morelos'
$ws.Cells.Item(1051, 2).Value = 'Villa Sola De Vega'
$ws.Cells.Item(1052, 2).Value = 'Villa De Etla'
$ws.Cells.Item(1053, 2).Value = 'Villa De Tamazulápam Del Progreso'
$ws.Cells.Item(1054, 2).Value = 'Villa De Tututepec'
$ws.Cells.Item(1055, 2).Value = 'Villa De Tututepec De Melchor Ocampo'
$ws.Cells.Item(1056, 2).Value = 'Villa De Zaachila'
$ws.Cells.Item(1057, 2).Value = 'Yutanduchi De Guerrero'
$ws.Cells.Item(1059, 2).Value = 'Zimatlán De Álvarez'
$ws.Cells.Item(1080, 2).Value = 'Chalchicomula De Sesma'
$ws.Cells.Item(1088, 2).Value = 'Chila De La Sal'
$ws.Cells.Item(1097, 2).Value = 'Cuetzalan Del Progreso'
$ws.Cells.Item(1110, 2).Value = 'Huehuetlán El Chico'
$ws.Cells.Item(1111, 2).Value = 'Huehuetlán El Grande'
$ws.Cells.Item(1117, 2).Value = 'Izúcar De Matamoros'
$ws.Cells.Item(1126, 2).Value = 'Los Reyes De Juárez'
$ws.Cells.Item(1132, 2).Value = 'Palmar De Bravo'
$ws.Cells.Item(1149, 2).Value = 'San Salvador El Seco'
$ws.Cells.Item(1154, 2).Value = 'Tecali De Herrera'
$ws.Cells.Item(1160, 2).Value = 'Tepanco De López'
$ws.Cells.Item(1161, 2).Value = 'Tepango De Rodríguez'
$ws.Cells.Item(1162, 2).Value = 'Tepatlaxco De Hidalgo'
$ws.Cells.Item(1165, 2).Value = 'Tepexi De Rodríguez'
$ws.Cells.Item(1167, 2).Value = 'Tetela De Ocampo'
$ws.Cells.Item(1172, 2).Value = 'Tlacotepec De Benito Juárez'
$ws.Cells.Item(1182, 2).Value = 'Tuzamapan De Galeana'
$ws.Cells.Item(1198, 2).Value = 'Amealco De Bonfil'
$ws.Cells.Item(1200, 2).Value = 'Cadereyta De Montes'
$ws.Cells.Item(1205, 2).Value = 'Jalpan De Serra'
$ws.Cells.Item(1206, 2).Value = 'Landa De Matamoros'
$ws.Cells.Item(1209, 2).Value = 'Pinal De Amoles'
$ws.Cells.Item(1212, 2).Value = 'San Juan Del Río'
$ws.Cells.Item(1225, 2).Value = 'Armadillo De Los Infante'
$ws.Cells.Item(1226, 2).Value = 'Axtla De Terrazas'
$ws.Cells.Item(1233, 2).Value = 'Ciudad Del Maíz'
$ws.Cells.Item(1243, 2).Value = 'Mexquitic De Carmona'
$ws.Cells.Item(1248, 2).Value = 'San Ciro De Acosta'
$ws.Cells.Item(1254, 2).Value = 'Santa María Del Río'
$ws.Cells.Item(1256, 2).Value = 'Soledad De Graciano Sánchez'
$ws.Cells.Item(1263, 2).Value = 'Tanquián De Escobedo'
$ws.Cells.Item(1268, 2).Value = 'Villa De Arista'
$ws.Cells.Item(1269, 2).Value = 'Villa De Arriaga'
$ws.Cells.Item(1270, 2).Value = 'Villa De Guadalupe'
$ws.Cells.Item(1271, 2).Value = 'Villa De Ramos'
$ws.Cells.Item(1272, 2).Value = 'Villa De Reyes'
$ws.Cells.Item(1273, 2).Value = 'Villa De La Paz'
$ws.Cells.Item(1295, 2).Value = 'Nacozari De García'
$ws.Cells.Item(1309, 2).Value = 'Jalpa De Méndez'
$ws.Cells.Item(1339, 2).Value = 'Soto La Marina'
$ws.Cells.Item(1350, 2).Value = 'Contla De Juan Cuamatzi'
$ws.Cells.Item(1355, 2).Value = 'Nanacamilpa De Mariano Arista'
$ws.Cells.Item(1357, 2).Value = 'Papalotla De Xicohténcatl'
$ws.Cells.Item(1359, 2).Value = 'San Pablo Del Monte'
$ws.Cells.Item(1360, 2).Value = 'Sanctórum De Lázaro Cárdenas'
$ws.Cells.Item(1363, 2).Value = 'Tepetitla De Lardizábal'
$ws.Cells.Item(1365, 2).Value = 'Tetla De La Solidaridad'
$ws.Cells.Item(1374, 1).Value = 'Veracruz De Ignacio De La Llave'
$ws.Cells.Item(1381, 2).Value = 'Alto Lucero De Gutiérrez Barrios'
$ws.Cells.Item(1385, 2).Value = 'Amatlán De Los Reyes'
$ws.Cells.Item(1399, 2).Value = 'Castillo De Teayo'
$ws.Cells.Item(1401, 2).Value = 'Cazones De Herrera'
$ws.Cells.Item(1419, 2).Value = 'Cosamaloapan De Carpio'
$ws.Cells.Item(1420, 2).Value = 'Cosautlán De Carvajal'
$ws.Cells.Item(1438, 2).Value = 'Hueyapan De Ocampo'
$ws.Cells.Item(1439, 2).Value = 'Ignacio De La Llave'
$ws.Cells.Item(1442, 2).Value = 'Ixhuacán De Los Reyes'
$ws.Cells.Item(1443, 2).Value = 'Ixhuatlán De Madero'
$ws.Cells.Item(1444, 2).Value = 'Ixhuatlán Del Café'
$ws.Cells.Item(1445, 2).Value = 'Ixhuatlán Del Sureste'
$ws.Cells.Item(1453, 2).Value = 'Juchique De Ferrer'
$ws.Cells.Item(1459, 2).Value = 'Las Vigas De Ramírez'
$ws.Cells.Item(1460, 2).Value = 'Lerdo De Tejada'
$ws.Cells.Item(1463, 2).Value = 'Martínez De La Torre'
$ws.Cells.Item(1469, 2).Value = 'Mixtla De Altamirano'
$ws.Cells.Item(1471, 2).Value = 'Nanchital De Lázaro Cárdenas Del Río'
$ws.Cells.Item(1480, 2).Value = 'Ozuluama De Mascareñas'
$ws.Cells.Item(1483, 2).Value = 'Paso De Ovejas'
$ws.Cells.Item(1484, 2).Value = 'Paso Del Macho'
$ws.Cells.Item(1488, 2).Value = 'Poza Rica De Hidalgo'
$ws.Cells.Item(1497, 2).Value = 'Sayula De Alemán'
$ws.Cells.Item(1500, 2).Value = 'Soledad De Doblado'
$ws.Cells.Item(1507, 2).Value = 'Tatahuicapan De Juárez'
$ws.Cells.Item(1539, 2).Value = 'Vega De Alatorre'
$ws.Cells.Item(1550, 2).Value = 'Zontecomatlán De López Y Fuentes'
$ws.Cells.Item(1551, 2).Value = 'Zozocolco De Hidalgo'
$ws.Cells.Item(1566, 2).Value = 'Cañitas De Felipe Pescador'
$ws.Cells.Item(1580, 2).Value = 'Mezquital Del Oro'
$ws.Cells.Item(1584, 2).Value = 'Nochistlán De Mejía'
$ws.Cells.Item(1585, 2).Value = 'Noria De Ángeles'
$ws.Cells.Item(1595, 2).Value = 'Teúl De González Ortega'
$ws.Cells.Item(1596, 2).Value = 'Tlaltenango De Sánchez Román'
$ws.Cells.Item(1602, 2).Value = 'Villa De Cos'
$ws.Cells.Item(1606, 1).Value = 'Total'

# Remove the trailing footer/metadata rows (1608-1612)
$ws.Rows('1608:1612').Delete()

